$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.719.93'
$ws.Range("E2").Value = '  -1.49%  '
$ws.Range("D3").Value = '3.403.78'
$ws.Range("E3").Value = '  +3.40%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '254.86'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '648.54'
$ws.Range("E6").Value = '  +3.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.45'
$ws.Range("E7").Value = '  -0.26%  '
$ws.Range("E8").Value = '  +5.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.05'
$ws.Range("E9").Value = '  +6.61%  '
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").Value = '3.400.30'
$ws.Range("E11").Value = '  +3.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.211'
$ws.Range("E12").Value = '  +4.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.21'
$ws.Range("E13").Value = '  +2.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.24'
$ws.Range("E14").Value = '  +13.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000256'
$ws.Range("E15").Value = '  +2.51%  '
$ws.Range("D16").Value = '97.315.40'
$ws.Range("E16").Value = '  -1.55%  '
$ws.Range("D17").Value = '4.039.73'
$ws.Range("E17").Value = '  +3.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.47'
$ws.Range("E18").Value = '  +32.25%  '
$ws.Range("D19").Value = '3.421.86'
$ws.Range("E19").Value = '  +3.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.31'
$ws.Range("E20").Value = '  +13.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.70'
$ws.Range("E21").Value = '  +14.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.487'
$ws.Range("E22").Value = '  +40.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.40'
$ws.Range("E23").Value = '  -1.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '500.22'
$ws.Range("E24").Value = '  +1.99%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.05'
$ws.Range("E26").Value = '  +6.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.81'
$ws.Range("E27").Value = '  +7.35%  '
$ws.Range("E28").Value = '  +3.81%  '
$ws.Range("D29").Value = '3.580.02'
$ws.Range("E29").Value = '  +3.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.152'
$ws.Range("E30").Value = '  +9.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.199'
$ws.Range("E31").Value = '  +5.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.998'
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("E33").Value = '  +5.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.564'
$ws.Range("E35").Value = '  +17.57%  '
$ws.Range("E36").Value = '  +5.05%  '
$ws.Range("E37").Value = '  +14.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.64'
$ws.Range("E38").Value = '  +4.72%  '
$ws.Range("E39").Value = '  +1.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.40'
$ws.Range("E40").Value = '  +12.88%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '509.18'
$ws.Range("E41").Value = '  +3.75%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.854'
$ws.Range("E43").Value = '  +10.24%  '
$ws.Range("E44").Value = '  -3.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0412'
$ws.Range("E45").Value = '  +21.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.48'
$ws.Range("E46").Value = '  +14.83%  '
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.20'
$ws.Range("E48").Value = '  +2.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.13'
$ws.Range("E49").Value = '  +10.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.56'
$ws.Range("E50").Value = '  +13.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.41'
$ws.Range("E51").Value = '  +10.84%  '
